$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Nppa"
$ws.Range("C2").Value = "Npr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2077686666666667
$ws.Range("H2").Value = 0.623306
$ws.Range("I2").Value = 0.3307450952508051
$ws.Range("J2").Value = 0.3307450952508051
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04090066666666667
$ws.Range("N2").Value = 0.122702
$ws.Range("O2").Value = 0.09164231251535751
$ws.Range("P2").Value = 0.0916423125153575
$ws.Range("Q2").Value = 0.008497876979111113
$ws.Range("R2").Value = 0.076480892812
$ws.Range("S2").Value = 0.03031024538189597
$ws.Range("T2").Value = 0.03031024538189596

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Nppa"
$ws.Range("C3").Value = "Npr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2077686666666667
$ws.Range("H3").Value = 0.623306
$ws.Range("I3").Value = 0.3307450952508051
$ws.Range("J3").Value = 0.3307450952508051
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.3809676666666666
$ws.Range("N3").Value = 1.142903
$ws.Range("O3").Value = 0.8535987506376395
$ws.Range("P3").Value = 0.8535987506376393
$ws.Range("Q3").Value = 0.07915314414644445
$ws.Range("R3").Value = 0.7123782973180001
$ws.Range("S3").Value = 0.2823236000856143
$ws.Range("T3").Value = 0.2823236000856142

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Nppa"
$ws.Range("C4").Value = "Npr3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2077686666666667
$ws.Range("H4").Value = 0.623306
$ws.Range("I4").Value = 0.3307450952508051
$ws.Range("J4").Value = 0.3307450952508051
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02443933333333333
$ws.Range("N4").Value = 0.07331799999999999
$ws.Range("O4").Value = 0.05475893684700315
$ws.Range("P4").Value = 0.05475893684700314
$ws.Range("Q4").Value = 0.005077727700888888
$ws.Range("R4").Value = 0.045699549308
$ws.Range("S4").Value = 0.01811124978329488
$ws.Range("T4").Value = 0.01811124978329488

$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("B5").Value = "Nppa"
$ws.Range("C5").Value = "Npr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.420415
$ws.Range("H5").Value = 1.261245
$ws.Range("I5").Value = 0.6692549047491948
$ws.Range("J5").Value = 0.6692549047491948
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04090066666666667
$ws.Range("N5").Value = 0.122702
$ws.Range("O5").Value = 0.09164231251535751
$ws.Range("P5").Value = 0.0916423125153575
$ws.Range("Q5").Value = 0.01719525377666667
$ws.Range("R5").Value = 0.15475728399
$ws.Range("S5").Value = 0.06133206713346153
$ws.Range("T5").Value = 0.06133206713346152

$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Nppa"
$ws.Range("C6").Value = "Npr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.420415
$ws.Range("H6").Value = 1.261245
$ws.Range("I6").Value = 0.6692549047491948
$ws.Range("J6").Value = 0.6692549047491948
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3809676666666666
$ws.Range("N6").Value = 1.142903
$ws.Range("O6").Value = 0.8535987506376395
$ws.Range("P6").Value = 0.8535987506376393
$ws.Range("Q6").Value = 0.1601645215816667
$ws.Range("R6").Value = 1.441480694235
$ws.Range("S6").Value = 0.571275150552025
$ws.Range("T6").Value = 0.571275150552025

$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Nppa"
$ws.Range("C7").Value = "Npr3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.420415
$ws.Range("H7").Value = 1.261245
$ws.Range("I7").Value = 0.6692549047491948
$ws.Range("J7").Value = 0.6692549047491948
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.02443933333333333
$ws.Range("N7").Value = 0.07331799999999999
$ws.Range("O7").Value = 0.05475893684700315
$ws.Range("P7").Value = 0.05475893684700314
$ws.Range("Q7").Value = 0.01027466232333333
$ws.Range("R7").Value = 0.09247196090999998
$ws.Range("S7").Value = 0.03664768706370827
$ws.Range("T7").Value = 0.03664768706370826
